# v0.8 - Almost there!
#
# 1) Relocate the "_GoBack" bookmark from inside the
#    "...spends the majority of its' [time idle]" paragraph to a brand new,
#    otherwise-empty paragraph placed right after the title
#    ("Event-driven Programming") and before the "Definition" heading.
#    The two runs that used to straddle the bookmark ("spends the majority
#    of its'" / " time idle") become a single merged run once the bookmark
#    is gone.
#
# 2) In the "Fortunately for today's computer users..." paragraph, merge
#    the two runs into one, and move the <w:lastRenderedPageBreak/> marker
#    (which used to sit on the second run) to the front of the single
#    merged run.

$d = $word.ActiveDocument
$wordMlNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------
# Change 1: move the _GoBack bookmark into its own new empty paragraph
# ---------------------------------------------------------------------

# Remove the bookmark from its current location (this does not delete any
# text, only the bookmark markers themselves).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# Create a brand new, completely empty paragraph right after the first
# paragraph (the "Event-driven Programming" heading).
$titlePara = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.InsertXML('<w:p xmlns:w="' + $wordMlNs + '"></w:p>')

# Put the bookmark markers inside that new empty paragraph.
$newPara = $d.Paragraphs.Item(2)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newParaRange.InsertXML('<w:p xmlns:w="' + $wordMlNs + '"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# Merge the two runs that used to be split by the bookmark into a single
# run: "spends the majority of its<rsquo>" + " time idle".
$mergeRange = $d.Content
$quote = [char]0x2019
$null = $mergeRange.Find.Execute("spends the majority of its" + $quote + " time idle", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergedText = $mergeRange.Text
$mergeRange.Delete()
$mergeRange.InsertAfter($mergedText)

# ---------------------------------------------------------------------
# Change 2: merge the "Fortunately..." paragraph's two runs into one,
# keeping <w:lastRenderedPageBreak/> as the first child of the run.
# ---------------------------------------------------------------------

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Fortunately for today*") {
        $targetPara = $candidate
        break
    }
}

$paraTextRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
$fullText = $paraTextRange.Text
$escapedText = $fullText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
$replacementXml = '<w:p xmlns:w="' + $wordMlNs + '"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $escapedText + '</w:t></w:r></w:p>'
$paraTextRange.InsertXML($replacementXml)
